$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3078.739
$ws.Range("J64").Value = 3053.25
$ws.Range("L64").Value = 3053.25
$ws.Range("N64").Value = -3549.25

$ws.Range("H67").Value = 3078.739
$ws.Range("J67").Value = 3053.25
$ws.Range("L67").Value = 3053.25
$ws.Range("N67").Value = -4769.25

$ws.Range("H80").Value = 1735.1818
$ws.Range("I80").Value = 1883.8572
$ws.Range("J80").Value = 1475
$ws.Range("K80").Value = 5651.571599999999
$ws.Range("L80").Value = 4425
$ws.Range("M80").Value = -4653.571599999999
$ws.Range("N80").Value = -6421

$ws.Range("H82").Value = 7656.8335
$ws.Range("I82").Value = 800
$ws.Range("J82").Value = 9028.200000000001
$ws.Range("K82").Value = 2400
$ws.Range("L82").Value = 27084.6
$ws.Range("M82").Value = -1994
$ws.Range("N82").Value = -27896.6

$ws.Range("H83").Value = 1735.1818
$ws.Range("I83").Value = 1883.8572
$ws.Range("J83").Value = 1475
$ws.Range("K83").Value = 16954.7148
$ws.Range("L83").Value = 13275
$ws.Range("M83").Value = -11962.7148
$ws.Range("N83").Value = -23259

$ws.Range("H85").Value = 7656.8335
$ws.Range("I85").Value = 800
$ws.Range("J85").Value = 9028.200000000001
$ws.Range("K85").Value = 2400
$ws.Range("L85").Value = 27084.6
$ws.Range("M85").Value = -996
$ws.Range("N85").Value = -29892.6

$ws.Range("H87").Value = 30270.8
$ws.Range("J87").Value = 30270.8
$ws.Range("L87").Value = 30270.8
$ws.Range("N87").Value = -32766.8

$ws.Range("H90").Value = 30270.8
$ws.Range("J90").Value = 30270.8
$ws.Range("L90").Value = 90812.39999999999
$ws.Range("N90").Value = -103292.4

$ws.Range("H100").Value = 4749.8335
$ws.Range("I100").Value = 4500
$ws.Range("J100").Value = 4999.6665
$ws.Range("K100").Value = 4500
$ws.Range("L100").Value = 4999.6665
$ws.Range("M100").Value = -3959
$ws.Range("N100").Value = -6081.6665

$ws.Range("H103").Value = 3756692.8
$ws.Range("I103").Value = 7512698
$ws.Range("J103").Value = 687.5
$ws.Range("K103").Value = 22538094
$ws.Range("L103").Value = 2062.5
$ws.Range("M103").Value = -22537508
$ws.Range("N103").Value = -3234.5

$ws.Range("H137").Value = 4350283
$ws.Range("I137").Value = 1516.8334
$ws.Range("J137").Value = 9094391
$ws.Range("K137").Value = 4550.5002
$ws.Range("L137").Value = 27283173
$ws.Range("M137").Value = -2000.5002
$ws.Range("N137").Value = -27288273

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 933.8333
$ws.Range("I20").Value = 941.41174
$ws.Range("J20").Value = 915.4286
$ws.Range("K20").Value = 941.41174
$ws.Range("L20").Value = 915.4286
$ws.Range("M20").Value = -694.41174
$ws.Range("N20").Value = -1409.4286

$ws.Range("H86").Value = 20001930
$ws.Range("I86").Value = 1961.0526
$ws.Range("J86").Value = 83335170
$ws.Range("K86").Value = 1961.0526
$ws.Range("L86").Value = 83335170
$ws.Range("M86").Value = -838.0526
$ws.Range("N86").Value = -83337416

$ws.Range("H89").Value = 20001930
$ws.Range("I89").Value = 1961.0526
$ws.Range("J89").Value = 83335170
$ws.Range("K89").Value = 9805.262999999999
$ws.Range("L89").Value = 416675850
$ws.Range("M89").Value = -4189.262999999999
$ws.Range("N89").Value = -416687082

$ws.Range("H105").Value = 3786.8823
$ws.Range("I105").Value = 2660
$ws.Range("J105").Value = 4028.3572
$ws.Range("K105").Value = 2660
$ws.Range("L105").Value = 4028.3572
$ws.Range("M105").Value = -913
$ws.Range("N105").Value = -7522.3572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 138.8
$ws.Range("I7").Value = 123
$ws.Range("J7").Value = 202
$ws.Range("K7").Value = 123
$ws.Range("L7").Value = 202
$ws.Range("M7").Value = -10
$ws.Range("N7").Value = -428

$ws.Range("H31").Value = 5053749.5
$ws.Range("I31").Value = 4727.5947
$ws.Range("J31").Value = 11495605
$ws.Range("K31").Value = 4727.5947
$ws.Range("L31").Value = 11495605
$ws.Range("M31").Value = -4432.5947
$ws.Range("N31").Value = -11496195

$ws.Range("H34").Value = 5053749.5
$ws.Range("I34").Value = 4727.5947
$ws.Range("J34").Value = 11495605
$ws.Range("K34").Value = 4727.5947
$ws.Range("L34").Value = 11495605
$ws.Range("M34").Value = -4525.5947
$ws.Range("N34").Value = -11496009

$ws.Range("H62").Value = 2301.25
$ws.Range("I62").Value = 2270.7693
$ws.Range("K62").Value = 2270.7693
$ws.Range("M62").Value = -1646.7693

$ws.Range("H65").Value = 2301.25
$ws.Range("I65").Value = 2270.7693
$ws.Range("K65").Value = 11353.8465
$ws.Range("M65").Value = -8233.8465

$ws.Range("H86").Value = 2370.8
$ws.Range("I86").Value = 2750
$ws.Range("K86").Value = 2750
$ws.Range("M86").Value = -1627

$ws.Range("H89").Value = 2370.8
$ws.Range("I89").Value = 2750
$ws.Range("K89").Value = 13750
$ws.Range("M89").Value = -8134

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 854.44183
$ws.Range("I68").Value = 585.9091
$ws.Range("J68").Value = 1740.6
$ws.Range("K68").Value = 1757.7273
$ws.Range("L68").Value = 5221.799999999999
$ws.Range("M68").Value = -946.7273
$ws.Range("N68").Value = -6843.799999999999

$ws.Range("H71").Value = 854.44183
$ws.Range("I71").Value = 585.9091
$ws.Range("J71").Value = 1740.6
$ws.Range("K71").Value = 5273.1819
$ws.Range("L71").Value = 15665.4
$ws.Range("M71").Value = -1217.1819
$ws.Range("N71").Value = -23777.4

$ws.Range("H107").Value = 930.7308
$ws.Range("I107").Value = 420.96295
$ws.Range("J107").Value = 1481.28
$ws.Range("K107").Value = 1262.88885
$ws.Range("L107").Value = 4443.84
$ws.Range("M107").Value = 657.1111500000002
$ws.Range("N107").Value = -8283.84

$ws.Range("H121").Value = 1207.1538
$ws.Range("I121").Value = 261.14285
$ws.Range("J121").Value = 2310.8333
$ws.Range("K121").Value = 783.4285500000001
$ws.Range("L121").Value = 6932.499899999999
$ws.Range("M121").Value = 526.5714499999999
$ws.Range("N121").Value = -9552.499899999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 44731.6
$ws.Range("J70").Value = 5251
$ws.Range("L70").Value = 5251
$ws.Range("N70").Value = -5791

$ws.Range("H73").Value = 44731.6
$ws.Range("J73").Value = 5251
$ws.Range("L73").Value = 5251
$ws.Range("N73").Value = -7123

$ws.Range("H122").Value = 3510039.8
$ws.Range("I122").Value = 6061655
$ws.Range("J122").Value = 1568.75
$ws.Range("K122").Value = 18184965
$ws.Range("L122").Value = 4706.25
$ws.Range("M122").Value = -18182515
$ws.Range("N122").Value = -9606.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5749.2915
$ws.Range("I40").Value = 7617.909
$ws.Range("J40").Value = 4168.154
$ws.Range("K40").Value = 7617.909
$ws.Range("L40").Value = 4168.154
$ws.Range("M40").Value = -7481.909
$ws.Range("N40").Value = -4440.154
